$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
$ws.Activate()

# Task "Implement fix" (row 6) moved from "In progress" to "Done",
# and its effort on Day 6 (column K) bumped from 1 to 2.
$ws.Range("F6").Value = "Done"
$ws.Range("K6").Value = 2

# Reflect the user's last selection on the sheet.
$ws.Range("F6").Select() | Out-Null
